# readxl "types.xlsx" test fixture: add coverage for xls BoolErr-style
# error values (#DIV/0!, #N/A, #VALUE!) by inserting a new "errors" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before D. This shifts the old "date" column (D) to E
# and the old trailing integer column (E) to F - formulas referencing them
# are adjusted automatically by Excel.
$ws.Range("D1").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("D1").Value = "errors"

# Row 2: flip the boolean, add a formula-driven #DIV/0! error.
$ws.Range("C2").Value = $true
$ws.Range("D2").Formula = "=1/0"

# Row 3: flip the boolean, add a literal (non-formula) #N/A error value -
# this mirrors the xls BoolErr record case that isn't a formula result.
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = "#N/A"

# Row 4: add a formula-driven #VALUE! error (number + text).
$ws.Range("D4").Formula = "=F2+F4"

# Match the new active selection recorded in the sheet view.
$ws.Range("E3").Select() | Out-Null
